# Apply the "fetching coins" update to the workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet1: A1 price value update ---
$wsSheet1 = $wb.Worksheets.Item("Sheet1")
$wsSheet1.Range("A1").Value = "342,18"

# --- data: coin labels / sheet refs shift one column right (B<-removed, etc.) ---
$wsData = $wb.Worksheets.Item("data")
$wsData.Range("B1").Value = "ETHEREUM"
$wsData.Range("C1").Value = "LITECOIN"
$wsData.Range("D1").Value = "-"

$wsData.Range("B2").Value = "Sheet1"
$wsData.Range("D2").Value = ""

$wsData.Range("B3").Value = "A1"
$wsData.Range("D3").Value = ""

# --- Arkusz3: A1 price value update ---
# Plain assignment of "2,632075" would be mis-parsed as the number
# 2632075 (comma read as a thousands separator), which also drags in a
# new "Text" number-format style. Stage the literal text in a scratch
# cell explicitly formatted as Text, copy only its *value* onto A1 (so
# A1 keeps its original General-format style), then discard the scratch
# cell completely so no stray cell/style survives in the sheet.
$wsArkusz3 = $wb.Worksheets.Item("Arkusz3")
$scratch = $wsArkusz3.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "2,632075"
$scratch.Copy()
$wsArkusz3.Range("A1").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
$scratch.Clear()
